$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Update the count from 36 to 40 (exact, case-sensitive match keeps the
#    surrounding run structure/formatting untouched).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("36", $false, $false, $false, $false, $false, $true, 1, $false, "40", 2)

# ---------------------------------------------------------------------------
# 2) Append the new sentences (with proofing marks) to the paragraph that now
#    reads "...выполнены первые 40 пункта плана."
# ---------------------------------------------------------------------------
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "пункта плана\.") {
        $target = $para
    }
}

$r = $target.Range
$r.MoveEnd(1, -1)          # exclude the paragraph mark from the range
$insPos = $r.End           # position right before the paragraph mark

$rPr = '<w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:color w:val="0D1822"/><w:spacing w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'

$body = ""
$body += "<w:r>$rPr<w:t xml:space=`"preserve`"> </w:t></w:r>"
$body += "<w:r>$rPr<w:t xml:space=`"preserve`">Давай продолжать разработку, следуй строго по плану. Если </w:t></w:r>"
$body += '<w:proofErr w:type="gramStart"/>'
$body += "<w:r>$rPr<w:t>какие то</w:t></w:r>"
$body += '<w:proofErr w:type="gramEnd"/>'
$body += "<w:r>$rPr<w:t xml:space=`"preserve`"> пункты кажутся тебе слишком объемными для домохозяйки, разбивай на </w:t></w:r>"
$body += '<w:proofErr w:type="spellStart"/>'
$body += "<w:r>$rPr<w:t>подшаги</w:t></w:r>"
$body += '<w:proofErr w:type="spellEnd"/>'

$paraPr = '<w:pPr><w:rPr><w:rFonts w:ascii="Roboto" w:hAnsi="Roboto"/><w:color w:val="0D1822"/><w:spacing w:val="2"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body>' +
    '<w:p w14:paraId="559CA70E" w14:textId="271DA4E9" w:rsidR="00554335" w:rsidRDefault="00554335" w:rsidP="00554335">' +
    $paraPr + $body +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insRange = $d.Range($insPos, $insPos)
$insRange.InsertXML($xml)

# The InsertXML call above created a brand-new paragraph right after the
# target paragraph's mark; delete that now-redundant paragraph mark so the
# newly inserted runs become part of the original paragraph again (the
# surviving paragraph keeps the pPr/paraId we supplied above, which is
# identical to the original, so nothing observable changes there).
$markRange = $d.Range($insPos, $insPos + 1)
$markRange.Delete()
